# Rename worksheet tabs (sheetId order: 1=GNG, 2=NB, 3=RS, 4=TOL, 5=vSAT)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "GNG_TO-16512556085069742"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556126235564"
$wb.Worksheets.Item(3).Name = "RS_TO-16512556126255586"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651255612672559"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512556127515578"

# Sheet 1: GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255608473977.csv"
$ws1.Range("B3").Value = "GNG_stims-16512556084889762.csv"
$ws1.Range("B4").Value = "go_stims-16512556084909754.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556085049748.csv"

# Sheet 2: NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_1-16512556088239748.csv"
$ws2.Range("B3").Value = "OB-1651255609860827.csv"
$ws2.Range("B4").Value = "TB-1651255612611557.csv"
$ws2.Range("B5").Value = "TB-16512556114305809.csv"
$ws2.Range("B6").Value = "ZB-match_0-1651255609281974.csv"
$ws2.Range("B7").Value = "ZB-match_1-16512556087829742.csv"
$ws2.Range("B8").Value = "OB-1651255609585828.csv"
$ws2.Range("B9").Value = "OB-1651255609462827.csv"
$ws2.Range("B10").Value = "TB-16512556125685568.csv"

# Sheet 3: RS (swap eyes closed / eyes open)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4: TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556126395566.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556126285589.csv"
$ws4.Range("B4").Value = "MM_stims-16512556126555617.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556126405585.csv"
$ws4.Range("B6").Value = "MM_stims-16512556126715574.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255612657558.csv"

# Sheet 5: vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512556127195597.csv"
$ws5.Range("B3").Value = "SAT_stims-16512556126785576.csv"
$ws5.Range("B4").Value = "SAT_stims-16512556127035596.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512556127365596.csv"
